$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Colonne parfois vide" (C) / "Repas" (D) columns, added while starting
# the importfromexcel.js script that uploads this sheet to MongoDB.
# Column C first (it's only sometimes filled in), then column D.
$ws.Range("C2").Value = "Colonne parfois vide"
$ws.Range("C4").Value = "Rempli"

$ws.Range("D2").Value = "Repas"
$ws.Range("D3").Value = "Thon"
$ws.Range("D4").Value = "Thon"
$ws.Range("D5").Value = "Viande"

# Match the author's final selection in the saved workbook.
$ws.Range("C15").Select()
